# Edit job_history-style workbook:
# - remove the outlier-affected "BP only" rows' distinctive formatting, promote
#   the next rows' number-cell style down the table
# - append a new data row (ukb51139_subset.csv, 26474 x 1081, "no events, remove
#   outliers") and a fresh trailing blank row
# - widen column D to fit the new, longer "Row Subset" text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Stash a copy of the "s=7/s=8" (K/N/O theme-font, bordered) number format
#    in a scratch range before we overwrite every row that currently carries
#    it (rows 23-26). We'll paste it back onto rows 27, 28, 30 and 31 later.
# ---------------------------------------------------------------------------
$ws.Range("K23:O23").Copy()
$ws.Range("Q201:U201").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 1. Rows 23-25: K/N go from the "s=7" style back to the plain "s=5" style
#    (copy format from row 22, which keeps that style throughout), and O goes
#    from "s=8" to "s=6".
# ---------------------------------------------------------------------------
$ws.Range("K22:O22").Copy()
foreach ($r in 23..25) {
    $ws.Range("K$r:O$r").PasteSpecial(-4122)
}

# Row 26: same style change, plus its row height grows to 19.5
$ws.Range("K26:O26").PasteSpecial(-4122)
$ws.Rows.Item(26).RowHeight = 19.5

# ---------------------------------------------------------------------------
# 2. Row 27: Column Subset changes from "BP only" to "all"; K/N/O pick up the
#    stashed "s=7/s=8" style (now vacated by rows 23-26); row height -> 19.5
# ---------------------------------------------------------------------------
$ws.Cells.Item(27, 3).Value = "all"
$ws.Range("Q201:U201").Copy()
$ws.Range("K27:O27").PasteSpecial(-4122)
$ws.Rows.Item(27).RowHeight = 19.5

# ---------------------------------------------------------------------------
# 3. Row 28 was blank; it now holds a new run (outliers removed). Row height
#    is unchanged.
# ---------------------------------------------------------------------------
$ws.Cells.Item(28, 1).Value = "ukb51139_subset.csv"
$ws.Cells.Item(28, 2).Value = "26474 x 1081"
$ws.Cells.Item(28, 3).Value = "all"
$ws.Cells.Item(28, 4).Value = "no events, remove outliers"
$ws.Cells.Item(28, 5).Value = "> 140/80"
$ws.Cells.Item(28, 6).Value = "zscore"
$ws.Cells.Item(28, 7).Value = "median"
$ws.Cells.Item(28, 8).Value = "none"
$ws.Cells.Item(28, 9).Value = 50
$ws.Cells.Item(28, 11).Value = 85
$ws.Cells.Item(28, 12).Value = "91.7 & 90.4"
$ws.Cells.Item(28, 13).Value = "79.5 & 75.2"
$ws.Cells.Item(28, 14).Value = 18
$ws.Cells.Item(28, 15).Value = 15.6

# I28 picks up the plain "s=5" style (same source as rows 23-25 above)
$ws.Range("I22").Copy()
$ws.Range("I28").PasteSpecial(-4122)
# K28/N28/O28 pick up the stashed "s=7/s=8" style
$ws.Range("Q201:U201").Copy()
$ws.Range("K28:O28").PasteSpecial(-4122)
# restore the values that PasteSpecial(formats) left untouched (format-only
# paste does not disturb values, but keep this explicit/idempotent)
$ws.Cells.Item(28, 11).Value = 85
$ws.Cells.Item(28, 14).Value = 18
$ws.Cells.Item(28, 15).Value = 15.6

# Row 29 is left exactly as-is (still blank).

# ---------------------------------------------------------------------------
# 4. Row 30: still blank, but I/K/N/O move from the borderless "s=9/s=10"
#    blank style to the bordered "s=7/s=8" style.
# ---------------------------------------------------------------------------
$ws.Range("Q201:U201").Copy()
$ws.Range("I30").PasteSpecial(-4122)
$ws.Range("K30:O30").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5. Row 31: still blank; I/K/N/O move from "s=5/s=6" to "s=7/s=8"; row
#    height shrinks back to 18.75.
# ---------------------------------------------------------------------------
$ws.Range("Q201:U201").Copy()
$ws.Range("I31").PasteSpecial(-4122)
$ws.Range("K31:O31").PasteSpecial(-4122)
$ws.Rows.Item(31).RowHeight = 18.75

# Row 32 is left exactly as-is (still blank).

# ---------------------------------------------------------------------------
# 6. Append a brand-new trailing blank row 33, matching row 32's formatting.
# ---------------------------------------------------------------------------
$ws.Range("A32:O32").Copy()
$ws.Range("A33:O33").PasteSpecial(-4122)
$ws.Rows.Item(33).RowHeight = 19.5

# ---------------------------------------------------------------------------
# 7. Clean up the scratch range used to stash formats.
# ---------------------------------------------------------------------------
$ws.Range("Q201:U201").Clear()

# ---------------------------------------------------------------------------
# 8. Widen column D to fit the new, longer text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 22.5
